$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text formatted so numeric-looking strings
# like "47.317.55" are not reinterpreted by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.317.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.497.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.86"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.91"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.499.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.859"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.256.92"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.97%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "250.74"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.24"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.137"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.64%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.27"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0795"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +5.86%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.45"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.978.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +11.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.33%  "
